$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
@("ECs","App","Fpr2","ECs",3,1,89.38217433333334,268.146523,0.2143552015363441,0.2175965347165783,1,0.3333333333333333,0.050657,0.151971,0.002123644810648064,0.002123644810648064,4.527832805203667,40.750495246833,0.000455214311378077,0.0004620977517658628),
@("ECs","App","Fpr2","FAPs",3,1,89.38217433333334,268.146523,0.2143552015363441,0.2175965347165783,3,1,1.006255666666667,3.018767,0.0421842909114609,0.0421842909114609,89.94131942190457,809.471874797141,0.009042422179993971,0.009179155521809942),
@("ECs","App","Fpr2","Inflammatory-Mac",3,1,89.38217433333334,268.146523,0.2143552015363441,0.2175965347165783,3,1,18.813815,56.441445,0.7887135162611822,0.7887135162611822,1681.619692205082,15134.57722984574,0.1690648447326043,0.1716213280225609),
@("ECs","App","Fpr2","Resolving-Mac",3,1,89.38217433333334,268.146523,0.2143552015363441,0.2175965347165783,3,1,3.983073,11.949219,0.1669785480167087,0.1669785480167088,356.015725268393,3204.141527415537,0.03579272031236772,0.03633395342044161),
@("FAPs","App","Fpr2","ECs",3,1,122.3539896666667,367.061969,0.2934277926151677,0.2978648075949286,1,0.3333333333333333,0.050657,0.151971,0.002123644810648064,0.002123644810648064,6.198086054544333,55.78277449089899,0.0006231364090871171,0.0006325590529236541),
@("FAPs","App","Fpr2","FAPs",3,1,122.3539896666667,367.061969,0.2934277926151677,0.2978648075949286,3,1,1.006255666666667,3.018767,0.0421842909114609,0.0421842909114609,123.1193954413581,1108.074558972223,0.01237804336518605,0.01256521569587079),
@("FAPs","App","Fpr2","Inflammatory-Mac",3,1,122.3539896666667,367.061969,0.2934277926151677,0.2978648075949286,3,1,18.813815,56.441445,0.7887135162611822,0.7887135162611822,2301.945326100578,20717.5079349052,0.2314304660822659,0.2349299997686566),
@("FAPs","App","Fpr2","Resolving-Mac",3,1,122.3539896666667,367.061969,0.2934277926151677,0.2978648075949286,3,1,3.983073,11.949219,0.1669785480167087,0.1669785480167088,487.344872683579,4386.10385415221,0.04899614675862863,0.0497370330774775),
@("Inflammatory-Mac","App","Fpr2","ECs",3,1,90.33462533333334,271.003876,0.2166393574945233,0.2199152301234996,1,0.3333333333333333,0.050657,0.151971,0.002123644810648064,0.002123644810648064,4.576081115510667,41.184730039596,0.0004600650473253751,0.0004670218372342447),
@("Inflammatory-Mac","App","Fpr2","FAPs",3,1,90.33462533333334,271.003876,0.2166393574945233,0.2199152301234996,3,1,1.006255666666667,3.018767,0.0421842909114609,0.0421842909114609,90.89972863787689,818.0975577408919,0.009138777679420946,0.009276968043390578),
@("Inflammatory-Mac","App","Fpr2","Inflammatory-Mac",3,1,90.33462533333334,271.003876,0.2166393574945233,0.2199152301234996,3,1,18.813815,56.441445,0.7887135162611822,0.7887135162611822,1699.538929115647,15295.85036204082,0.1708663894100687,0.1734501144300925),
@("Inflammatory-Mac","App","Fpr2","Resolving-Mac",3,1,90.33462533333334,271.003876,0.2166393574945233,0.2199152301234996,3,1,3.983073,11.949219,0.1669785480167087,0.1669785480167088,359.8094071303159,3238.284664172844,0.03617412535770818,0.03672112581278234),
@("MuSCs","App","Fpr2","ECs",2,1,18.634161,37.268322,0.0446882095496985,0.03024263611988591,1,0.3333333333333333,0.050657,0.151971,0.002123644810648064,0.002123644810648064,0.9439506937769999,5.663704162661999,0.00009490188430737047,0.00006422461725631342),
@("MuSCs","App","Fpr2","FAPs",2,1,18.634161,37.268322,0.0446882095496985,0.03024263611988591,3,1,1.006255666666667,3.018767,0.0421842909114609,0.0421842909114609,18.750730099829,112.504380598974,0.001885140431956807,0.001275764160010723),
@("MuSCs","App","Fpr2","Inflammatory-Mac",2,1,18.634161,37.268322,0.0446882095496985,0.03024263611988591,3,1,18.813815,56.441445,0.7887135162611822,0.7887135162611822,350.579657734215,2103.47794640529,0.03524619488935925,0.02385277587512265),
@("MuSCs","App","Fpr2","Resolving-Mac",2,1,18.634161,37.268322,0.0446882095496985,0.03024263611988591,3,1,3.983073,11.949219,0.1669785480167087,0.1669785480167088,74.22122355675299,445.3273413405179,0.007461972344075074,0.005049871467496221),
@("Resolving-Mac","App","Fpr2","ECs",3,1,96.27664699999998,288.829941,0.2308894388042666,0.2343807914451077,1,0.3333333333333333,0.050657,0.151971,0.002123644810648064,0.002123644810648064,4.877086107078999,43.893774963711,0.0004903271585501245,0.000497741551467989),
@("Resolving-Mac","App","Fpr2","FAPs",3,1,96.27664699999998,288.829941,0.2308894388042666,0.2343807914451077,3,1,1.006255666666667,3.018767,0.0421842909114609,0.0421842909114609,96.87892161141632,871.9102945027469,0.009739907254903131,0.009887187490378866),
@("Resolving-Mac","App","Fpr2","Inflammatory-Mac",3,1,96.27664699999998,288.829941,0.2308894388042666,0.2343807914451077,3,1,18.813815,56.441445,0.7887135162611822,0.7887135162611822,1811.331025478305,16301.97922930474,0.1821056211468841,0.1848592981647497),
@("Resolving-Mac","App","Fpr2","Resolving-Mac",3,1,96.27664699999998,288.829941,0.2308894388042666,0.2343807914451077,3,1,3.983073,11.949219,0.1669785480167087,0.1669785480167088,383.4769131962309,3451.292218766078,0.03855358324392916,0.03913656423851111)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Count; $j++) {
    $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
  }
}
